$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely (d, Televisor, 1, 47, TRUE), shifting rows 4 and 5 up.
$ws.Rows.Item(3).Delete()
